$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.252.47'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.498.32'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.17'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.08'
$ws.Range("E6").Value = '  +3.21%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.124'
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.24'
$ws.Range("E10").Value = '  +1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.097.14'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.120'
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000181'
$ws.Range("E14").Value = '  +3.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.502.37'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.91'
$ws.Range("E16").Value = '  -4.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.284.99'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.90'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("E19").Value = '  +2.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.61'
$ws.Range("E20").Value = '  -3.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '393.36'
$ws.Range("E21").Value = '  +3.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.571'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.641.42'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("E26").Value = '  +1.58%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.40'
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.27'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.520.77'
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.151'
$ws.Range("E33").Value = '  +5.12%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.43'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.90'
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '165.51'
$ws.Range("E39").Value = '  +3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0783'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.806'
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.29'
$ws.Range("E43").Value = '  -3.43%  '
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("E46").Value = '  +3.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.457.42'
$ws.Range("E47").Value = '  +1.56%  '
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.895'
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.216'
$ws.Range("E51").Value = '  +0.38%  '
